$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168419003486633
$ws.Range("B1").Value = 2.167053937911987
$ws.Range("C1").Value = 10.38233280181885
$ws.Range("D1").Value = 2.55781888961792
$ws.Range("E1").Value = 1.255747199058533
